# The "process_fcs" transformation (flow_cyt) and its corresponding entry
# in the Views sheet (flow_cy) were removed from the test fixture.

$wb = $excel.ActiveWorkbook

# "Transformations" sheet: remove the row defining flow_cyt = process_fcs(...)
$wsTransformations = $wb.Worksheets.Item("Transformations")
$wsTransformations.Rows.Item(2).Delete()
$wsTransformations.Range("I18").Select()

# "Views" sheet: remove the row mapping flow_cy -> flow_cyt
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Rows.Item(2).Delete()
$wsViews.Range("C11").Select()
